$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$rows = @(45, 46, 47, 50, 51, 52, 53, 54, 55, 56, 57)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = $false
}
